# Insert two new data rows (901:902) into the "Plátano" sheet, pushing the
# existing rows 901..969 down to 903..971, and populate the new rows with a
# fresh Pintón / Primera Pintón price pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 901 downward (by two) to make room for the new pair.
$ws.Rows("901:902").Insert()

# --- New row 901: Pintón ---
$ws.Range("A901").Value = 7
$ws.Range("B901").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C901").Value = "Ñuble"
$ws.Range("D901").Value = 45106
$ws.Range("E901").Value = 16
$ws.Range("F901").Value = "Fruta"
$ws.Range("G901").Value = 100108
$ws.Range("H901").Value = "Tropicales y subtropicales"
$ws.Range("I901").Value = 100108006
$ws.Range("J901").Value = "Plátano"
$ws.Range("K901").Value = "Sin especificar"
$ws.Range("L901").Value = "Pintón"
$ws.Range("M901").Value = 120
$ws.Range("N901").Value = 14000
$ws.Range("O901").Value = 14000
$ws.Range("P901").Value = 14000
$ws.Range("Q901").Value = "$/caja 20 kilos"
$ws.Range("R901").Value = "Ecuador"
$ws.Range("S901").Value = 700
$ws.Range("T901").Value = 20

# --- New row 902: Primera Pintón ---
$ws.Range("A902").Value = 7
$ws.Range("B902").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C902").Value = "Ñuble"
$ws.Range("D902").Value = 45106
$ws.Range("E902").Value = 16
$ws.Range("F902").Value = "Fruta"
$ws.Range("G902").Value = 100108
$ws.Range("H902").Value = "Tropicales y subtropicales"
$ws.Range("I902").Value = 100108006
$ws.Range("J902").Value = "Plátano"
$ws.Range("K902").Value = "Sin especificar"
$ws.Range("L902").Value = "Primera Pintón"
$ws.Range("M902").Value = 150
$ws.Range("N902").Value = 14500
$ws.Range("O902").Value = 14500
$ws.Range("P902").Value = 14500
$ws.Range("Q902").Value = "$/caja 20 kilos"
$ws.Range("R902").Value = "Ecuador"
$ws.Range("S902").Value = 725
$ws.Range("T902").Value = 20
